# FAST_holdings.xlsx update:
#  - Bump the "as of" date in the confidential disclosure text from
#    2021-03-31 to 2021-04-05
#  - Refresh the Weight / Percent Change values for rows 2-10 (D2:E10)
#
# The sheet is protected, so it must be unprotected before the cells can
# be written to, then re-protected afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# --- Update the disclosure paragraph's date -----------------------------
$disclosure = $ws.Range("A13").Value2
$disclosure = $disclosure -replace "2021-03-31", "2021-04-05"
$ws.Range("A13").Value2 = $disclosure

# --- Update Weight (D) / Percent Change (E) figures ----------------------
$ws.Range("D2").Value2 = 0.1054902413638055
$ws.Range("E2").Value2 = 0.01430007820355272

$ws.Range("D3").Value2 = 0.1057853718931584
$ws.Range("E3").Value2 = 0.01419741753870629

$ws.Range("D4").Value2 = 0.1164608639431013
$ws.Range("E4").Value2 = 0.01234868795190502

$ws.Range("D5").Value2 = 0.136983059104259
$ws.Range("E5").Value2 = 0.00766315095583403

$ws.Range("D6").Value2 = 0.1330876017010799
$ws.Range("E6").Value2 = 0.01176133103843946

$ws.Range("D7").Value2 = 0.1443258004408463
$ws.Range("E7").Value2 = 0.01323543618974021

$ws.Range("D8").Value2 = 0.1293461831793628
$ws.Range("E8").Value2 = 0.01948051948051943

$ws.Range("D9").Value2 = 0.1285208783743866
$ws.Range("E9").Value2 = 0.01552075012167298

$ws.Range("D10").Value2 = 0.9999999999999999
$ws.Range("E10").Value2 = 0.01348823206602878

$ws.Protect()
